$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("no")

# Update header row text (order matters for shared-string table ordering)
$ws.Range("B1").Value = "Company name "
$ws.Range("C1").Value = "Record ID"
$ws.Range("A1").Value = "Website URL"

# Best-fit the data columns (matches the author's column-width autofit)
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

# Update selection to match diff (I8)
$ws.Range("I8").Select()
